$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'306.78"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'-3.24%"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'40.99"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'-2.13%"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'5.048"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'-2.97%"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'0.07616"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'-5.72%"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'4.245"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'-2.92%"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'1.596"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'-8.90%"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'0.9055"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'-2.64%"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.09963"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'-11.16%"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'0.1768"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'-4.89%"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.09216"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'-1.13%"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.04418"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'-3.31%"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = "'-0.11%"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'0.001257"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'-2.88%"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'0.005821"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'-1.56%"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'3.365"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'0.28%"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "'-3.62%"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'-2.92%"
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'6.771"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'-8.34%"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = "'-3.14%"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "'11.52%"
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'0.04162"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'-0.10%"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'0.001215"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'-2.24%"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'0.004064"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'-6.01%"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'0.0001301"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'6.31%"
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'0.0003009"
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'0.02416"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'-6.27%"
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.05144"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'-5.57%"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.007831"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'-2.63%"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'-6.25%"
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'0.007074"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'-6.55%"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.001949"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'-6.54%"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.008258"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'-0.18%"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.3042"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'-3.20%"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'0.00006385"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'-5.83%"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "'-0.21%"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "'-26.94%"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'0.006077"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'79.13%"
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'0.00002101"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'-0.21%"
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'0.0002001"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'-0.21%"
$c.Style = "Normal"
